# Update financial ratio figures on the "ESS" worksheet.
# Only columns D:G for rows 13-17 and 28-29 change, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESS")

# Row 13 - Gross Margin
$ws.Range("D13").Value = 0.6916
$ws.Range("E13").Value = 0.6976
$ws.Range("F13").Value = 0.7059
$ws.Range("G13").Value = 0.7051

# Row 14 - EBIT Margin
$ws.Range("D14").Value = 0.3289
$ws.Range("E14").Value = 0.3265
$ws.Range("F14").Value = 0.3317
$ws.Range("G14").Value = 0.3295

# Row 15 - EBT margin
$ws.Range("D15").Value = 0.4205
$ws.Range("E15").Value = 0.4394
$ws.Range("F15").Value = 0.4462
$ws.Range("G15").Value = 0.3181

# Row 16 - Net Profit Margin
$ws.Range("D16").Value = 0.4
$ws.Range("E16").Value = 0.4181
$ws.Range("F16").Value = 0.4247
$ws.Range("G16").Value = 0.3008

# Row 17 - Free Cash Flow Margin
$ws.Range("D17").Value = 0.2884
$ws.Range("E17").Value = 0.2332
$ws.Range("F17").Value = 0.0618
$ws.Range("G17").Value = 0.3275

# Row 28 - EBITDA Margin
$ws.Range("D28").Value = 0.6651
$ws.Range("E28").Value = 0.6549
$ws.Range("F28").Value = 0.6502
$ws.Range("G28").Value = 0.6452

# Row 29 - Operating Cash Flow Margin
$ws.Range("D29").Value = 0.5983
$ws.Range("E29").Value = 0.6248
$ws.Range("F29").Value = 0.6264
$ws.Range("G29").Value = 0.6294
